$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6000.625
$ws.Range("I62").Value = 2252.5
$ws.Range("J62").Value = 7250
$ws.Range("K62").Value = 2252.5
$ws.Range("L62").Value = 7250
$ws.Range("M62").Value = -1628.5
$ws.Range("N62").Value = -8498

$ws.Range("H65").Value = 6000.625
$ws.Range("I65").Value = 2252.5
$ws.Range("J65").Value = 7250
$ws.Range("K65").Value = 11262.5
$ws.Range("L65").Value = 36250
$ws.Range("M65").Value = -8142.5
$ws.Range("N65").Value = -42490

$ws.Range("H132").Value = 1985.8
$ws.Range("I132").Value = 1913.4642
$ws.Range("K132").Value = 5740.392599999999
$ws.Range("M132").Value = -3210.392599999999

$ws.Range("H137").Value = 11332.523
$ws.Range("I137").Value = 4602.533
$ws.Range("J137").Value = 15071.407
$ws.Range("K137").Value = 13807.599
$ws.Range("L137").Value = 45214.221
$ws.Range("M137").Value = -11257.599
$ws.Range("N137").Value = -50314.221

$ws.Range("H138").Value = 6486.933
$ws.Range("I138").Value = 8648.333000000001
$ws.Range("J138").Value = 6154.41
$ws.Range("K138").Value = 25944.999
$ws.Range("L138").Value = 18463.23
$ws.Range("M138").Value = -20804.999
$ws.Range("N138").Value = -28743.23

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 198.33333
$ws.Range("I4").Value = 198.33333
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 198.33333
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -82.33332999999999
$ws.Range("N4").ClearContents()

$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H23").Value = 1158
$ws.Range("J23").Value = 1158
$ws.Range("L23").Value = 1158
$ws.Range("N23").Value = -1676

$ws.Range("H32").Value = 5057656.5
$ws.Range("I32").Value = 5439297
$ws.Range("K32").Value = 5439297
$ws.Range("M32").Value = -5439010

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H45").Value = 2240.9524
$ws.Range("I45").Value = 1898.8572
$ws.Range("J45").Value = 2412
$ws.Range("K45").Value = 1898.8572
$ws.Range("L45").Value = 2412
$ws.Range("M45").Value = -1521.8572
$ws.Range("N45").Value = -3166

$ws.Range("H61").Value = 22777208
$ws.Range("I61").Value = 33336142
$ws.Range("K61").Value = 33336142
$ws.Range("M61").Value = -33335930

$ws.Range("H74").Value = 5324829
$ws.Range("I74").Value = 7577368
$ws.Range("K74").Value = 7577368
$ws.Range("M74").Value = -7576494

$ws.Range("H77").Value = 5324829
$ws.Range("I77").Value = 7577368
$ws.Range("K77").Value = 37886840
$ws.Range("M77").Value = -37882472

$ws.Range("H136").Value = 22777208
$ws.Range("I136").Value = 33336142
$ws.Range("K136").Value = 100008426
$ws.Range("M136").Value = -100005876

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H19").Value = 3500
$ws.Range("I19").Value = 3000
$ws.Range("K19").Value = 3000
$ws.Range("M19").Value = -2827

$ws.Range("H22").Value = 899.8333
$ws.Range("I22").Value = 979.8
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 979.8
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -806.8
$ws.Range("N22").Value = -846

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H134").Value = 28209.121
$ws.Range("I134").Value = 3460.1724
$ws.Range("J134").Value = 88019.086
$ws.Range("K134").Value = 10380.5172
$ws.Range("L134").Value = 264057.258
$ws.Range("M134").Value = -7845.5172
$ws.Range("N134").Value = -269127.258

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -888
$ws.Range("N4").ClearContents()

$ws.Range("H7").Value = 1574.7097
$ws.Range("I7").Value = 203.70589
$ws.Range("J7").Value = 3239.5
$ws.Range("K7").Value = 203.70589
$ws.Range("L7").Value = 3239.5
$ws.Range("M7").Value = -90.70589000000001
$ws.Range("N7").Value = -3465.5

$ws.Range("H31").Value = 482074.3
$ws.Range("I31").Value = 11196.706
$ws.Range("J31").Value = 732228
$ws.Range("K31").Value = 11196.706
$ws.Range("L31").Value = 732228
$ws.Range("M31").Value = -10901.706
$ws.Range("N31").Value = -732818

$ws.Range("H34").Value = 482074.3
$ws.Range("I34").Value = 11196.706
$ws.Range("J34").Value = 732228
$ws.Range("K34").Value = 11196.706
$ws.Range("L34").Value = 732228
$ws.Range("M34").Value = -10994.706
$ws.Range("N34").Value = -732632

$ws.Range("H122").Value = 3481.75
$ws.Range("I122").Value = 2058.4707
$ws.Range("J122").Value = 5681.364
$ws.Range("K122").Value = 6175.4121
$ws.Range("L122").Value = 17044.092
$ws.Range("M122").Value = -3725.4121
$ws.Range("N122").Value = -21944.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("K14").Value = 1500
$ws.Range("M14").Value = -1327

$ws.Range("H68").Value = 3380.9
$ws.Range("I68").Value = 3081.818
$ws.Range("J68").Value = 3494.3447
$ws.Range("K68").Value = 9245.454000000002
$ws.Range("L68").Value = 10483.0341
$ws.Range("M68").Value = -8434.454000000002
$ws.Range("N68").Value = -12105.0341

$ws.Range("H71").Value = 3380.9
$ws.Range("I71").Value = 3081.818
$ws.Range("J71").Value = 3494.3447
$ws.Range("K71").Value = 27736.362
$ws.Range("L71").Value = 31449.1023
$ws.Range("M71").Value = -23680.362
$ws.Range("N71").Value = -39561.1023

$ws.Range("H97").Value = 1786747.1
$ws.Range("I97").Value = 5102552
$ws.Range("J97").Value = 1313.7693
$ws.Range("K97").Value = 15307656
$ws.Range("L97").Value = 3941.3079
$ws.Range("M97").Value = -15307160
$ws.Range("N97").Value = -4933.3079

$ws.Range("H127").Value = 52864.227
$ws.Range("J127").Value = 52864.227
$ws.Range("L127").Value = 158592.681
$ws.Range("N127").Value = -168512.681

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 218.18182
$ws.Range("I2").Value = 128.71428
$ws.Range("J2").Value = 374.75
$ws.Range("K2").Value = 128.71428
$ws.Range("L2").Value = 374.75
$ws.Range("M2").Value = -15.71428
$ws.Range("N2").Value = -600.75

$ws.Range("H97").Value = 1704.5264
$ws.Range("I97").Value = 1467.875
$ws.Range("K97").Value = 1467.875
$ws.Range("M97").Value = -971.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705

$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893

$ws.Range("H46").Value = 2714.3171
$ws.Range("I46").Value = 2337.2104
$ws.Range("J46").Value = 3040
$ws.Range("K46").Value = 2337.2104
$ws.Range("L46").Value = 3040
$ws.Range("M46").Value = -2149.2104
$ws.Range("N46").Value = -3416

$ws.Range("H55").Value = 58823936
$ws.Range("I55").Value = 76923520
$ws.Range("J55").Value = 292
$ws.Range("K55").Value = 76923520
$ws.Range("L55").Value = 292
$ws.Range("M55").Value = -76923347
$ws.Range("N55").Value = -638

$ws.Range("H122").Value = 5676.1304
$ws.Range("I122").Value = 4752.8335
$ws.Range("K122").Value = 14258.5005
$ws.Range("M122").Value = -11808.5005

$ws.Range("H136").Value = 30013.13
$ws.Range("I136").Value = 4811.9644
$ws.Range("J136").Value = 69214.94500000001
$ws.Range("K136").Value = 14435.8932
$ws.Range("L136").Value = 207644.835
$ws.Range("M136").Value = -11885.8932
$ws.Range("N136").Value = -212744.835

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14292020
$ws.Range("I62").Value = 4944.75
$ws.Range("K62").Value = 4944.75
$ws.Range("M62").Value = -4320.75

$ws.Range("H65").Value = 14292020
$ws.Range("I65").Value = 4944.75
$ws.Range("K65").Value = 24723.75
$ws.Range("M65").Value = -21603.75

$ws.Range("H116").Value = 78982
$ws.Range("J116").Value = 78982
$ws.Range("L116").Value = 78982
$ws.Range("N116").Value = -88160

$ws.Range("H119").Value = 51653.4
$ws.Range("J119").Value = 51653.4
$ws.Range("L119").Value = 51653.4
$ws.Range("N119").Value = -61329.4

$ws.Range("H122").Value = 8315.637000000001
$ws.Range("J122").Value = 19233
$ws.Range("L122").Value = 57699
$ws.Range("N122").Value = -62599

$ws.Range("H125").Value = 90248
$ws.Range("J125").Value = 90248
$ws.Range("L125").Value = 90248
$ws.Range("N125").Value = -100088
